$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = [DateTime]"05/17/2021"
$ws.Cells.Item(2, 13).Value = 58
$ws.Cells.Item(2, 14).Value = 10000
$ws.Cells.Item(2, 15).Value = 10000
$ws.Cells.Item(2, 16).Value = 10000
$ws.Cells.Item(2, 19).Value = 1000
$ws.Cells.Item(3, 4).Value = [DateTime]"05/17/2021"
$ws.Cells.Item(3, 13).Value = 65
$ws.Cells.Item(3, 14).Value = 9000
$ws.Cells.Item(3, 15).Value = 9000
$ws.Cells.Item(3, 16).Value = 9000
$ws.Cells.Item(3, 19).Value = 900
$ws.Cells.Item(4, 4).Value = [DateTime]"05/17/2021"
$ws.Cells.Item(4, 12).Value = "Segunda"
$ws.Cells.Item(4, 13).Value = 60
$ws.Cells.Item(4, 14).Value = 8000
$ws.Cells.Item(4, 15).Value = 8000
$ws.Cells.Item(4, 16).Value = 8000
$ws.Cells.Item(4, 19).Value = 800
$ws.Cells.Item(5, 4).Value = [DateTime]"05/18/2022"
$ws.Cells.Item(5, 12).Value = "Especial"
$ws.Cells.Item(5, 13).Value = 56
$ws.Cells.Item(5, 14).Value = 12000
$ws.Cells.Item(5, 15).Value = 12000
$ws.Cells.Item(5, 16).Value = 12000
$ws.Cells.Item(5, 19).Value = 1200
$ws.Cells.Item(6, 4).Value = [DateTime]"05/18/2022"
$ws.Cells.Item(6, 13).Value = 60
$ws.Cells.Item(7, 4).Value = [DateTime]"04/20/2021"
$ws.Cells.Item(7, 14).Value = 10000
$ws.Cells.Item(7, 15).Value = 10000
$ws.Cells.Item(7, 16).Value = 10000
$ws.Cells.Item(7, 19).Value = 1000
$ws.Cells.Item(8, 4).Value = [DateTime]"04/29/2021"
$ws.Cells.Item(8, 12).Value = "Primera"
$ws.Cells.Item(8, 13).Value = 45
$ws.Cells.Item(8, 14).Value = 10000
$ws.Cells.Item(8, 15).Value = 10000
$ws.Cells.Item(8, 16).Value = 10000
$ws.Cells.Item(8, 19).Value = 1000
$ws.Cells.Item(9, 4).Value = [DateTime]"04/26/2021"
$ws.Cells.Item(9, 13).Value = 48
$ws.Cells.Item(9, 14).Value = 10000
$ws.Cells.Item(9, 15).Value = 10000
$ws.Cells.Item(9, 16).Value = 10000
$ws.Cells.Item(9, 19).Value = 1000
$ws.Cells.Item(10, 4).Value = [DateTime]"05/05/2021"
$ws.Cells.Item(10, 12).Value = "Primera"
$ws.Cells.Item(10, 14).Value = 9000
$ws.Cells.Item(10, 15).Value = 9000
$ws.Cells.Item(10, 16).Value = 9000
$ws.Cells.Item(10, 19).Value = 900
$ws.Cells.Item(11, 4).Value = [DateTime]"04/28/2021"
$ws.Cells.Item(11, 13).Value = 47
$ws.Cells.Item(12, 4).Value = [DateTime]"05/10/2021"
$ws.Cells.Item(12, 12).Value = "Primera"
$ws.Cells.Item(12, 13).Value = 65
$ws.Cells.Item(12, 14).Value = 10000
$ws.Cells.Item(12, 15).Value = 10000
$ws.Cells.Item(12, 16).Value = 10000
$ws.Cells.Item(12, 19).Value = 1000
$ws.Cells.Item(13, 4).Value = [DateTime]"05/10/2021"
$ws.Cells.Item(13, 12).Value = "Segunda"
$ws.Cells.Item(13, 13).Value = 67
$ws.Cells.Item(13, 14).Value = 8000
$ws.Cells.Item(13, 15).Value = 8000
$ws.Cells.Item(13, 16).Value = 8000
$ws.Cells.Item(13, 19).Value = 800
$ws.Cells.Item(14, 4).Value = [DateTime]"05/07/2021"
$ws.Cells.Item(14, 12).Value = "Primera"
$ws.Cells.Item(14, 13).Value = 60
$ws.Cells.Item(14, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(15, 4).Value = [DateTime]"05/07/2021"
$ws.Cells.Item(15, 12).Value = "Segunda"
$ws.Cells.Item(15, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(16, 4).Value = [DateTime]"05/13/2021"
$ws.Cells.Item(16, 12).Value = "Primera"
$ws.Cells.Item(16, 13).Value = 56
$ws.Cells.Item(16, 14).Value = 9000
$ws.Cells.Item(16, 15).Value = 9000
$ws.Cells.Item(16, 16).Value = 9000
$ws.Cells.Item(16, 19).Value = 900
$ws.Cells.Item(17, 4).Value = [DateTime]"05/13/2021"
$ws.Cells.Item(17, 12).Value = "Segunda"
$ws.Cells.Item(17, 13).Value = 50
$ws.Cells.Item(17, 14).Value = 8000
$ws.Cells.Item(17, 15).Value = 8000
$ws.Cells.Item(17, 16).Value = 8000
$ws.Cells.Item(17, 18).Value = "Región Metropolitana"
$ws.Cells.Item(17, 19).Value = 800
$ws.Cells.Item(18, 4).Value = [DateTime]"04/16/2021"
$ws.Cells.Item(18, 13).Value = 45
$ws.Cells.Item(19, 4).Value = [DateTime]"05/27/2021"
$ws.Cells.Item(19, 12).Value = "Especial"
$ws.Cells.Item(19, 13).Value = 47
$ws.Cells.Item(19, 18).Value = "Región Metropolitana"
$ws.Cells.Item(20, 4).Value = [DateTime]"05/27/2021"
$ws.Cells.Item(20, 12).Value = "Primera"
$ws.Cells.Item(20, 13).Value = 50
$ws.Cells.Item(20, 14).Value = 9000
$ws.Cells.Item(20, 15).Value = 9000
$ws.Cells.Item(20, 16).Value = 9000
$ws.Cells.Item(20, 18).Value = "Región Metropolitana"
$ws.Cells.Item(20, 19).Value = 900
$ws.Cells.Item(21, 4).Value = [DateTime]"05/27/2021"
$ws.Cells.Item(21, 12).Value = "Segunda"
$ws.Cells.Item(21, 13).Value = 58
$ws.Cells.Item(21, 14).Value = 8000
$ws.Cells.Item(21, 15).Value = 8000
$ws.Cells.Item(21, 16).Value = 8000
$ws.Cells.Item(21, 18).Value = "Región Metropolitana"
$ws.Cells.Item(21, 19).Value = 800
$ws.Cells.Item(22, 4).Value = [DateTime]"04/23/2021"
$ws.Cells.Item(23, 4).Value = [DateTime]"04/22/2021"
$ws.Cells.Item(23, 13).Value = 45
$ws.Cells.Item(23, 14).Value = 10000
$ws.Cells.Item(23, 15).Value = 10000
$ws.Cells.Item(23, 16).Value = 10000
$ws.Cells.Item(23, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(23, 19).Value = 1000
$ws.Cells.Item(24, 4).Value = [DateTime]"04/22/2021"
$ws.Cells.Item(24, 13).Value = 48
$ws.Cells.Item(24, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(25, 4).Value = [DateTime]"04/15/2021"
$ws.Cells.Item(25, 13).Value = 45
$ws.Cells.Item(26, 4).Value = [DateTime]"05/12/2021"
$ws.Cells.Item(26, 12).Value = "Primera"
$ws.Cells.Item(26, 13).Value = 45
$ws.Cells.Item(27, 4).Value = [DateTime]"05/12/2021"
$ws.Cells.Item(27, 12).Value = "Segunda"
$ws.Cells.Item(27, 13).Value = 48
$ws.Cells.Item(27, 14).Value = 7000
$ws.Cells.Item(27, 15).Value = 7000
$ws.Cells.Item(27, 16).Value = 7000
$ws.Cells.Item(27, 19).Value = 700
$ws.Cells.Item(28, 4).Value = [DateTime]"05/03/2021"
$ws.Cells.Item(28, 13).Value = 68
$ws.Cells.Item(28, 14).Value = 10000
$ws.Cells.Item(28, 15).Value = 10000
$ws.Cells.Item(28, 16).Value = 10000
$ws.Cells.Item(28, 19).Value = 1000
$ws.Cells.Item(29, 4).Value = [DateTime]"05/03/2021"
$ws.Cells.Item(29, 12).Value = "Segunda"
$ws.Cells.Item(29, 13).Value = 57
$ws.Cells.Item(29, 14).Value = 8000
$ws.Cells.Item(29, 15).Value = 8000
$ws.Cells.Item(29, 16).Value = 8000
$ws.Cells.Item(29, 19).Value = 800
$ws.Cells.Item(30, 12).Value = "Primera"
$ws.Cells.Item(30, 13).Value = 56
$ws.Cells.Item(30, 14).Value = 10000
$ws.Cells.Item(30, 15).Value = 10000
$ws.Cells.Item(30, 16).Value = 10000
$ws.Cells.Item(30, 19).Value = 1000
$ws.Cells.Item(31, 4).Value = [DateTime]"05/06/2021"
$ws.Cells.Item(31, 12).Value = "Segunda"
$ws.Cells.Item(31, 13).Value = 40
$ws.Cells.Item(31, 14).Value = 8000
$ws.Cells.Item(31, 15).Value = 8000
$ws.Cells.Item(31, 16).Value = 8000
$ws.Cells.Item(31, 19).Value = 800
$ws.Cells.Item(32, 4).Value = [DateTime]"04/21/2021"
$ws.Cells.Item(32, 12).Value = "Primera"
$ws.Cells.Item(32, 13).Value = 40
$ws.Cells.Item(32, 14).Value = 10000
$ws.Cells.Item(32, 15).Value = 10000
$ws.Cells.Item(32, 16).Value = 10000
$ws.Cells.Item(32, 19).Value = 1000
